$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the sheet and its matching defined name -------------------
# (sheet rename auto-updates the defined name's formula reference; the
# name's own identifier must be renamed separately)
$ws.Name = "blood_lead"
$wb.Names.Item("Blood_Lead_Level").Name = "blood_lead"

# --- Relabel the five "*_5yavg" headers to "*_c1115" -------------------
$ws.Range("G1").Value = "_ebll_c1115"
$ws.Range("M1").Value = "_w_ebll_c1115"
$ws.Range("S1").Value = "_b_ebll_c1115"
$ws.Range("Y1").Value = "_a_ebll_c1115"
$ws.Range("AE1").Value = "_o_ebll_c1115"

# --- Widen column A ------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 49

# --- Force portrait page orientation (adds <pageSetup .../>) ------------
$ws.PageSetup.Orientation = 1
